# Doc Change / Move Animation and Card Animation Change
# Update the CharacterList column's type descriptor from an int-ref to a
# string-ref, and change the referenced data from numeric player IDs to
# the new character names (Zhouzhou / Timbuktu).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# C2 holds the type descriptor for the "CharacterList" column.
$ws.Range("C2").Value = "list#sep=|,string#ref=TbPlayerInfo"

# C4 holds the actual data for the CharacterList column; swap the
# numeric player IDs for the new character names.
$ws.Range("C4").Value = "Zhouzhou|Timbuktu"

# Move the active selection to C5 (matches the saved cursor position).
$null = $ws.Range("C5").Select()
